$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "GUTIERREZ CARLOS TERESA DE JESUS"
$ws.Range("B2").Value = 149

$ws.Range("A3").Value = "FERNANDEZ VALDERAS ERNESTO ALI"
$ws.Range("B3").Value = 144

$ws.Range("A4").Value = "CONTRERAS VALDERRAMA JULIA ALEJANDRA"
$ws.Range("B4").Value = 126

$ws.Range("A5").Value = "VALLE MAGALLAN EDUAR"
$ws.Range("B5").Value = 126

$ws.Range("A6").Value = "ZAVALETA MANAY JORGE LUIS"
$ws.Range("B6").Value = 106

$ws.Range("A7").Value = "ROMERO CHANAME YOSSELY TRINIDAD"
$ws.Range("B7").Value = 100

$ws.Range("A8").Value = "CAMACHO LINARES JUDITH ARLETT"
$ws.Range("B8").Value = 99

$ws.Range("A9").Value = "HUMPIRE CASTILLO IRWIN DEIMER"
$ws.Range("B9").Value = 99

$ws.Range("A10").Value = "SEVERINO AVALOS MARJORIE ISABEL"
$ws.Range("B10").Value = 98

$ws.Range("A11").Value = "HIDALGO CUBAS LUISA YVONE"
$ws.Range("B11").Value = 98

$ws.Range("A12").Value = "BALLENA ESQUÉN ASTRID CAROLINA"
$ws.Range("B12").Value = 93

$ws.Range("A13").Value = "SENADOR ARBOLEDA GIANCARLOS EXEBIO"
$ws.Range("B13").Value = 91

$ws.Range("A14").Value = "ZEVALLOS PACHECO ZOILA XIMENA"
$ws.Range("B14").Value = 88
